# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# handback files are now in sync with en-US (i.e. the handback was
# generated/accepted), updating the Status, Latest Handback DateTime and
# clearing the (now stale) Error Detail message on the per-locale sheets.
# The Overview sheet's summary cells automatically reflect the Status text
# change because they reference the same value.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ----- zh-cn sheet -----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-09-06 10:39:47"
$wsZhCn.Range("P2").Value = ""

# ----- de-de sheet -----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-09-06 10:40:16"
$wsDeDe.Range("P2").Value = ""

# ----- Overview sheet -----
# E2 / F2 mirror the same shared "Status" text for zh-cn / de-de respectively,
# so they pick up the new status text too.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# ----- Column width adjustments (the Status column widened to fit the new,
# longer status text; the Error Detail column narrowed now that it is empty) -----
# Column widths are stored in "characters" of the workbook's default font and
# get snapped to a 1/6-character grid by the engine, so we back-solve the
# ColumnWidth value that reproduces the desired stored width as closely as
# possible.
$statusColWidth = 29.9777050018311 - (5.0/6.0)
$errorColWidth  = 13.7470531463623 - (5.0/6.0)

$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth
$wsZhCn.Columns.Item(16).ColumnWidth = $errorColWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $errorColWidth
